$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1255.9166
$ws.Cells.Item(6, 9).Value = 1619.125
$ws.Cells.Item(6, 11).Value = 4857.375
$ws.Cells.Item(6, 13).Value = -4745.375
$ws.Cells.Item(100, 8).Value = 3999.8147
$ws.Cells.Item(100, 9).Value = 3071.0715
$ws.Cells.Item(100, 11).Value = 3071.0715
$ws.Cells.Item(100, 13).Value = -2530.0715
$ws.Cells.Item(113, 8).Value = 7743.1113
$ws.Cells.Item(113, 9).Value = 7736.4
$ws.Cells.Item(113, 11).Value = 7736.4
$ws.Cells.Item(113, 13).Value = -4482.4
$ws.Cells.Item(137, 8).Value = 4632958.5
$ws.Cells.Item(137, 9).Value = 3061.5715
$ws.Cells.Item(137, 11).Value = 9184.7145
$ws.Cells.Item(137, 13).Value = -6634.7145
$ws.Cells.Item(138, 8).Value = 2351
$ws.Cells.Item(138, 9).Value = 1088.9
$ws.Cells.Item(138, 10).Value = 2666.525
$ws.Cells.Item(138, 11).Value = 3266.7
$ws.Cells.Item(138, 12).Value = 7999.575000000001
$ws.Cells.Item(138, 13).Value = 1873.3
$ws.Cells.Item(138, 14).Value = -18279.575
$ws.Cells.Item(141, 8).Value = 6246.975
$ws.Cells.Item(141, 9).Value = 5242.074
$ws.Cells.Item(141, 11).Value = 15726.222
$ws.Cells.Item(141, 13).Value = -10546.222

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3468.6667
$ws.Cells.Item(63, 9).Value = 1930
$ws.Cells.Item(63, 10).Value = 4567.7144
$ws.Cells.Item(63, 11).Value = 1930
$ws.Cells.Item(63, 12).Value = 4567.7144
$ws.Cells.Item(63, 13).Value = -1244
$ws.Cells.Item(63, 14).Value = -5939.7144
$ws.Cells.Item(66, 8).Value = 3468.6667
$ws.Cells.Item(66, 9).Value = 1930
$ws.Cells.Item(66, 10).Value = 4567.7144
$ws.Cells.Item(66, 11).Value = 9650
$ws.Cells.Item(66, 12).Value = 22838.572
$ws.Cells.Item(66, 13).Value = -6218
$ws.Cells.Item(66, 14).Value = -29702.572
$ws.Cells.Item(74, 8).Value = 3027.6086
$ws.Cells.Item(74, 9).Value = 3036.8
$ws.Cells.Item(74, 10).Value = 2966.3333
$ws.Cells.Item(74, 11).Value = 3036.8
$ws.Cells.Item(74, 12).Value = 2966.3333
$ws.Cells.Item(74, 13).Value = -2162.8
$ws.Cells.Item(74, 14).Value = -4714.3333
$ws.Cells.Item(77, 8).Value = 3027.6086
$ws.Cells.Item(77, 9).Value = 3036.8
$ws.Cells.Item(77, 10).Value = 2966.3333
$ws.Cells.Item(77, 11).Value = 15184
$ws.Cells.Item(77, 12).Value = 14831.6665
$ws.Cells.Item(77, 13).Value = -10816
$ws.Cells.Item(77, 14).Value = -23567.6665
$ws.Cells.Item(132, 8).Value = 4485.722
$ws.Cells.Item(132, 9).Value = 4294.1
$ws.Cells.Item(132, 10).Value = 4725.25
$ws.Cells.Item(132, 11).Value = 12882.3
$ws.Cells.Item(132, 12).Value = 14175.75
$ws.Cells.Item(132, 13).Value = -10352.3
$ws.Cells.Item(132, 14).Value = -19235.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(116, 8).Value = 114999
$ws.Cells.Item(116, 10).Value = 114999
$ws.Cells.Item(116, 12).Value = 114999
$ws.Cells.Item(116, 14).Value = -124177

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 4001.75
$ws.Cells.Item(10, 9).Value = 2669
$ws.Cells.Item(10, 10).Value = 8000
$ws.Cells.Item(10, 11).Value = 2669
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = -2530
$ws.Cells.Item(10, 14).Value = -8278
$ws.Cells.Item(31, 8).Value = 11534.692
$ws.Cells.Item(31, 9).Value = 4308
$ws.Cells.Item(31, 10).Value = 12477.305
$ws.Cells.Item(31, 11).Value = 4308
$ws.Cells.Item(31, 12).Value = 12477.305
$ws.Cells.Item(31, 13).Value = -4013
$ws.Cells.Item(31, 14).Value = -13067.305
$ws.Cells.Item(34, 8).Value = 11534.692
$ws.Cells.Item(34, 9).Value = 4308
$ws.Cells.Item(34, 10).Value = 12477.305
$ws.Cells.Item(34, 11).Value = 4308
$ws.Cells.Item(34, 12).Value = 12477.305
$ws.Cells.Item(34, 13).Value = -4106
$ws.Cells.Item(34, 14).Value = -12881.305
$ws.Cells.Item(102, 8).Value = 30725
$ws.Cells.Item(102, 10).Value = 30725
$ws.Cells.Item(102, 12).Value = 30725
$ws.Cells.Item(102, 14).Value = -35593
$ws.Cells.Item(132, 8).Value = 4418.3057
$ws.Cells.Item(132, 9).Value = 3694.5806
$ws.Cells.Item(132, 10).Value = 8905.4
$ws.Cells.Item(132, 11).Value = 11083.7418
$ws.Cells.Item(132, 12).Value = 26716.2
$ws.Cells.Item(132, 13).Value = -8553.7418
$ws.Cells.Item(132, 14).Value = -31776.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 99.5
$ws.Cells.Item(17, 10).Value = 99.5
$ws.Cells.Item(17, 12).Value = 298.5
$ws.Cells.Item(17, 14).Value = -636.5
$ws.Cells.Item(113, 8).Value = 3179.8
$ws.Cells.Item(113, 9).Value = 999.6667
$ws.Cells.Item(113, 10).Value = 4114.143
$ws.Cells.Item(113, 11).Value = 2999.0001
$ws.Cells.Item(113, 12).Value = 12342.429
$ws.Cells.Item(113, 13).Value = -829.0001000000002
$ws.Cells.Item(113, 14).Value = -16682.429

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 10052.05
$ws.Cells.Item(2, 10).Value = 28649
$ws.Cells.Item(2, 12).Value = 28649
$ws.Cells.Item(2, 14).Value = -28875
$ws.Cells.Item(102, 8).Value = 1943.1428
$ws.Cells.Item(102, 9).Value = 1804.1538
$ws.Cells.Item(102, 11).Value = 1804.1538
$ws.Cells.Item(102, 13).Value = -182.1538
$ws.Cells.Item(107, 8).Value = 1135.091
$ws.Cells.Item(107, 9).Value = 1089.4
$ws.Cells.Item(107, 11).Value = 1089.4
$ws.Cells.Item(107, 13).Value = 830.5999999999999
$ws.Cells.Item(122, 8).Value = 10003.5
$ws.Cells.Item(122, 9).Value = 9999
$ws.Cells.Item(122, 11).Value = 29997
$ws.Cells.Item(122, 13).Value = -27547
$ws.Cells.Item(126, 8).Value = 2541.3333
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 14).Value = -16940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 1080001.8
$ws.Cells.Item(40, 8).Value = 6750
$ws.Cells.Item(40, 10).Value = 6750
$ws.Cells.Item(40, 12).Value = 6750
$ws.Cells.Item(40, 14).Value = -7022
$ws.Cells.Item(122, 8).Value = 15595.034
$ws.Cells.Item(122, 9).Value = 17129.777
$ws.Cells.Item(122, 11).Value = 51389.33099999999
$ws.Cells.Item(122, 13).Value = -48939.33099999999
$ws.Cells.Item(132, 8).Value = 5699.7144
$ws.Cells.Item(132, 9).Value = 5400.6665
$ws.Cells.Item(132, 10).Value = 5924
$ws.Cells.Item(132, 11).Value = 16201.9995
$ws.Cells.Item(132, 12).Value = 17772
$ws.Cells.Item(132, 13).Value = -13671.9995
$ws.Cells.Item(132, 14).Value = -22832
$ws.Cells.Item(136, 8).Value = 7679.278
$ws.Cells.Item(136, 9).Value = 5657.1816
$ws.Cells.Item(136, 11).Value = 16971.5448
$ws.Cells.Item(136, 13).Value = -14421.5448

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 386079.62
$ws.Cells.Item(2, 9).Value = 1586.25
$ws.Cells.Item(2, 11).Value = 1586.25
$ws.Cells.Item(2, 13).Value = -1474.25
$ws.Cells.Item(62, 8).Value = 5024.75
$ws.Cells.Item(62, 9).Value = 3165
$ws.Cells.Item(62, 11).Value = 3165
$ws.Cells.Item(62, 13).Value = -2541
$ws.Cells.Item(65, 8).Value = 5024.75
$ws.Cells.Item(65, 9).Value = 3165
$ws.Cells.Item(65, 11).Value = 15825
$ws.Cells.Item(65, 13).Value = -12705
$ws.Cells.Item(132, 8).Value = 3538.1667
$ws.Cells.Item(132, 9).Value = 3211.25
$ws.Cells.Item(132, 11).Value = 9633.75
$ws.Cells.Item(132, 13).Value = -7103.75
